$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new row with the YouTube link (new shared string, row 16)
$ws.Range("A16").Value = "https://www.youtube.com/watch?v=WHZn2cJNOkc"

# Scroll the view down (best effort - mirrors Excel's topLeftCell="A8")
try { $excel.ActiveWindow.ScrollRow = 8 } catch {}

# Update the selection to the newly added cell
$ws.Range("A16").Select()
